$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.912.37"
$ws.Range("E2").Value = "  -0.16%  "

$ws.Range("D3").Value = "1.548.96"
$ws.Range("E3").Value = "  -0.03%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "206.62"
$ws.Range("E5").Value = "  +0.37%  "

$ws.Range("D6").Value = "0.486"
$ws.Range("E6").Value = "  +0.41%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").Value = "22.08"
$ws.Range("E8").Value = "  +2.96%  "

$ws.Range("E9").Value = "  -0.55%  "

$ws.Range("E10").Value = "  +0.73%  "

$ws.Range("D11").Value = "0.0856"
$ws.Range("E11").Value = "  -0.13%  "

$ws.Range("D12").Value = "1.770.19"
$ws.Range("E12").Value = "  +0.12%  "

$ws.Range("D13").Value = "1.550.50"
$ws.Range("E13").Value = "  +0.23%  "

$ws.Range("E14").Value = "  +0.97%  "

$ws.Range("E15").Value = "  +1.25%  "

$ws.Range("D16").Value = "26.921.81"
$ws.Range("E16").Value = "  -0.02%  "

$ws.Range("D17").Value = "61.63"
$ws.Range("E17").Value = "  +0.08%  "

$ws.Range("D18").Value = "217.55"
$ws.Range("E18").Value = "  +1.59%  "

$ws.Range("D19").Value = "0.0₃0696"
$ws.Range("E19").Value = "  +1.68%  "

$ws.Range("D20").Value = "7.27"
$ws.Range("E20").Value = "  +0.67%  "

$ws.Range("E21").Value = "  -0.12%  "

$ws.Range("E22").Value = "  +0.54%  "

$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("E24").Value = "  +0.90%  "

$ws.Range("D25").Value = "154.38"
$ws.Range("E25").Value = "  +1.09%  "

$ws.Range("D26").Value = "6.62"
$ws.Range("E26").Value = "  -0.45%  "

$ws.Range("D27").Value = "14.91"
$ws.Range("E27").Value = "  +0.52%  "

$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("E30").Value = "  +1.56%  "

$ws.Range("D31").Value = "1.08"
$ws.Range("E31").Value = "  -0.75%  "

$ws.Range("E32").Value = "  -0.26%  "

$ws.Range("D33").Value = "1.415.31"
$ws.Range("E33").Value = "  +3.20%  "

$ws.Range("D34").Value = "3.08"
$ws.Range("E34").Value = "  +4.11%  "

$ws.Range("E35").Value = "  +2.34%  "

$ws.Range("D36").Value = "0.970"
$ws.Range("E36").Value = "  +0.23%  "

$ws.Range("E37").Value = "  +0.22%  "

$ws.Range("E38").Value = "  +0.67%  "

$ws.Range("D39").Value = "0.521"
$ws.Range("E39").Value = "  +0.59%  "

$ws.Range("E40").Value = "  +0.11%  "

$ws.Range("D41").Value = "5.75"
$ws.Range("E41").Value = "  +5.14%  "

$ws.Range("E42").Value = "  -0.15%  "

$ws.Range("D43").Value = "2.33"
$ws.Range("E43").Value = "  +4.73%  "

$ws.Range("D44").Value = "0.992"
$ws.Range("E44").Value = "  +0.55%  "

$ws.Range("D45").Value = "64.37"
$ws.Range("E45").Value = "  +1.28%  "

$ws.Range("E46").Value = "  +0.36%  "

$ws.Range("D47").Value = "1.683.78"
$ws.Range("E47").Value = "  +0.17%  "

$ws.Range("D48").Value = "87.69"
$ws.Range("E48").Value = "  +1.70%  "

$ws.Range("D49").Value = "0.0519"
$ws.Range("E49").Value = "  +2.36%  "

$ws.Range("E50").Value = "  +5.56%  "

$ws.Range("E51").Value = "  -0.01%  "
